$wb = $excel.ActiveWorkbook

# ---- Sheet 1: LP1912 ----
$ws1 = $wb.Worksheets.Item("LP1912")
$ws1.Range('A2').Value = 'Última actualización: 06:52:34'
$ws1.Range('A3').Value = 'Total filas: 69'

$ws1.Cells.Item(44,1).Value = '06:52:34'
$ws1.Cells.Item(44,2).Value = '06:59'
$ws1.Cells.Item(44,3).Value = '14_ABASTO'
$ws1.Cells.Item(44,4).Value = 7
$ws1.Cells.Item(44,5).Value = 'LP1912'

$ws1.Cells.Item(48,1).Value = '06:52:34'
$ws1.Cells.Item(48,2).Value = '07:05'
$ws1.Cells.Item(48,3).Value = '15_ABASTO'
$ws1.Cells.Item(48,4).Value = 13
$ws1.Cells.Item(48,5).Value = 'LP1912'

$ws1.Cells.Item(49,1).Value = '06:52:34'
$ws1.Cells.Item(49,2).Value = '07:05'
$ws1.Cells.Item(49,3).Value = '23_HERNANDEZ'
$ws1.Cells.Item(49,4).Value = 13
$ws1.Cells.Item(49,5).Value = 'LP1912'

$ws1.Cells.Item(50,1).Value = '06:24:49'
$ws1.Cells.Item(50,2).Value = '07:06'
$ws1.Cells.Item(50,3).Value = '225_GOMEZ'
$ws1.Cells.Item(50,4).Value = 42
$ws1.Cells.Item(50,5).Value = 'LP1912'

$ws1.Cells.Item(51,1).Value = '06:52:34'
$ws1.Cells.Item(51,2).Value = '07:07'
$ws1.Cells.Item(51,3).Value = '225_GOMEZ'
$ws1.Cells.Item(51,4).Value = 15
$ws1.Cells.Item(51,5).Value = 'LP1912'

$ws1.Cells.Item(52,1).Value = '06:52:34'
$ws1.Cells.Item(52,2).Value = '07:11'
$ws1.Cells.Item(52,3).Value = '215A_EL PATO'
$ws1.Cells.Item(52,4).Value = 19
$ws1.Cells.Item(52,5).Value = 'LP1912'

$ws1.Cells.Item(53,1).Value = '05:23:05'
$ws1.Cells.Item(53,2).Value = '07:12'
$ws1.Cells.Item(53,3).Value = '215A_EL PATO'
$ws1.Cells.Item(53,4).Value = 109
$ws1.Cells.Item(53,5).Value = 'LP1912'

$ws1.Cells.Item(54,1).Value = '06:52:34'
$ws1.Cells.Item(54,2).Value = '07:15'
$ws1.Cells.Item(54,3).Value = '11_ETCHEVERRY'
$ws1.Cells.Item(54,4).Value = 23
$ws1.Cells.Item(54,5).Value = 'LP1912'

$ws1.Cells.Item(55,1).Value = '05:23:05'
$ws1.Cells.Item(55,2).Value = '07:16'
$ws1.Cells.Item(55,3).Value = '11_ETCHEVERRY'
$ws1.Cells.Item(55,4).Value = 113
$ws1.Cells.Item(55,5).Value = 'LP1912'

$ws1.Cells.Item(56,1).Value = '06:52:34'
$ws1.Cells.Item(56,2).Value = '07:16'
$ws1.Cells.Item(56,3).Value = '16_SANTA ANA'
$ws1.Cells.Item(56,4).Value = 24
$ws1.Cells.Item(56,5).Value = 'LP1912'

$ws1.Cells.Item(57,1).Value = '06:52:34'
$ws1.Cells.Item(57,2).Value = '07:21'
$ws1.Cells.Item(57,3).Value = '26_HERNANDEZ'
$ws1.Cells.Item(57,4).Value = 29
$ws1.Cells.Item(57,5).Value = 'LP1912'

$ws1.Cells.Item(58,1).Value = '06:52:34'
$ws1.Cells.Item(58,2).Value = '07:23'
$ws1.Cells.Item(58,3).Value = '10_OLMOS'
$ws1.Cells.Item(58,4).Value = 31
$ws1.Cells.Item(58,5).Value = 'LP1912'

$ws1.Cells.Item(59,1).Value = '06:52:34'
$ws1.Cells.Item(59,2).Value = '07:31'
$ws1.Cells.Item(59,3).Value = '11_ETCHEVERRY'
$ws1.Cells.Item(59,4).Value = 39
$ws1.Cells.Item(59,5).Value = 'LP1912'

$ws1.Cells.Item(60,1).Value = '06:52:34'
$ws1.Cells.Item(60,2).Value = '07:32'
$ws1.Cells.Item(60,3).Value = '84_COLONIA URQUIZA-ESC 49'
$ws1.Cells.Item(60,4).Value = 40
$ws1.Cells.Item(60,5).Value = 'LP1912'

$ws1.Cells.Item(61,1).Value = '06:52:34'
$ws1.Cells.Item(61,2).Value = '07:36'
$ws1.Cells.Item(61,3).Value = '27_EL RETIRO'
$ws1.Cells.Item(61,4).Value = 44
$ws1.Cells.Item(61,5).Value = 'LP1912'

$ws1.Cells.Item(62,1).Value = '06:52:34'
$ws1.Cells.Item(62,2).Value = '07:39'
$ws1.Cells.Item(62,3).Value = '10_OLMOS'
$ws1.Cells.Item(62,4).Value = 47
$ws1.Cells.Item(62,5).Value = 'LP1912'

$ws1.Cells.Item(63,1).Value = '05:54:50'
$ws1.Cells.Item(63,2).Value = '07:46'
$ws1.Cells.Item(63,3).Value = '16_SANTA ANA'
$ws1.Cells.Item(63,4).Value = 112
$ws1.Cells.Item(63,5).Value = 'LP1912'

$ws1.Cells.Item(64,1).Value = '06:52:34'
$ws1.Cells.Item(64,2).Value = '07:47'
$ws1.Cells.Item(64,3).Value = '14_ABASTO'
$ws1.Cells.Item(64,4).Value = 55
$ws1.Cells.Item(64,5).Value = 'LP1912'

$ws1.Cells.Item(65,1).Value = '06:52:34'
$ws1.Cells.Item(65,2).Value = '07:51'
$ws1.Cells.Item(65,3).Value = '215D_EL PATO'
$ws1.Cells.Item(65,4).Value = 59
$ws1.Cells.Item(65,5).Value = 'LP1912'

$ws1.Cells.Item(66,1).Value = '06:52:34'
$ws1.Cells.Item(66,2).Value = '07:58'
$ws1.Cells.Item(66,3).Value = '16_SANTA ANA'
$ws1.Cells.Item(66,4).Value = 66
$ws1.Cells.Item(66,5).Value = 'LP1912'

$ws1.Cells.Item(67,1).Value = '06:24:49'
$ws1.Cells.Item(67,2).Value = '08:05'
$ws1.Cells.Item(67,3).Value = '23_HERNANDEZ'
$ws1.Cells.Item(67,4).Value = 101
$ws1.Cells.Item(67,5).Value = 'LP1912'

$ws1.Cells.Item(68,1).Value = '06:52:34'
$ws1.Cells.Item(68,2).Value = '08:06'
$ws1.Cells.Item(68,3).Value = '23_HERNANDEZ'
$ws1.Cells.Item(68,4).Value = 74
$ws1.Cells.Item(68,5).Value = 'LP1912'

$ws1.Cells.Item(69,1).Value = '06:52:34'
$ws1.Cells.Item(69,2).Value = '08:12'
$ws1.Cells.Item(69,3).Value = '15_ABASTO'
$ws1.Cells.Item(69,4).Value = 80
$ws1.Cells.Item(69,5).Value = 'LP1912'

$ws1.Cells.Item(70,1).Value = '06:52:34'
$ws1.Cells.Item(70,2).Value = '08:21'
$ws1.Cells.Item(70,3).Value = '26_HERNANDEZ'
$ws1.Cells.Item(70,4).Value = 89
$ws1.Cells.Item(70,5).Value = 'LP1912'

$ws1.Cells.Item(71,1).Value = '06:52:34'
$ws1.Cells.Item(71,2).Value = '08:22'
$ws1.Cells.Item(71,3).Value = '16_P MOR-SANTA ANA'
$ws1.Cells.Item(71,4).Value = 90
$ws1.Cells.Item(71,5).Value = 'LP1912'

$ws1.Cells.Item(72,1).Value = '06:52:34'
$ws1.Cells.Item(72,2).Value = '08:23'
$ws1.Cells.Item(72,3).Value = '215B_EL PATO'
$ws1.Cells.Item(72,4).Value = 91
$ws1.Cells.Item(72,5).Value = 'LP1912'

$ws1.Cells.Item(73,1).Value = '06:52:34'
$ws1.Cells.Item(73,2).Value = '08:27'
$ws1.Cells.Item(73,3).Value = '84_COLONIA URQUIZA-ESC 49'
$ws1.Cells.Item(73,4).Value = 95
$ws1.Cells.Item(73,5).Value = 'LP1912'

$ws1.Cells.Item(74,1).Value = '06:52:34'
$ws1.Cells.Item(74,2).Value = '08:42'
$ws1.Cells.Item(74,3).Value = '81_EL PELIGRO'
$ws1.Cells.Item(74,4).Value = 110
$ws1.Cells.Item(74,5).Value = 'LP1912'

# ---- Sheet 2: LP1912-215 ----
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Range('A2').Value = 'Última actualización: 06:52:34'

$ws2.Cells.Item(19,1).Value = '06:52:34'
$ws2.Cells.Item(19,2).Value = '07:11'
$ws2.Cells.Item(19,3).Value = '215A_EL PATO'
$ws2.Cells.Item(19,4).Value = 19
$ws2.Cells.Item(19,5).Value = 'LP1912'

$ws2.Cells.Item(21,1).Value = '06:52:34'
$ws2.Cells.Item(21,2).Value = '07:51'
$ws2.Cells.Item(21,3).Value = '215D_EL PATO'
$ws2.Cells.Item(21,4).Value = 59
$ws2.Cells.Item(21,5).Value = 'LP1912'

$ws2.Cells.Item(22,1).Value = '06:52:34'
$ws2.Cells.Item(22,2).Value = '08:23'
$ws2.Cells.Item(22,3).Value = '215B_EL PATO'
$ws2.Cells.Item(22,4).Value = 91
$ws2.Cells.Item(22,5).Value = 'LP1912'

# ---- Sheet 3: 6203-6173 ----
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Range('A2').Value = 'Última actualización: 06:52:34'
$ws3.Range('A3').Value = 'Total filas: 11'

$ws3.Cells.Item(13,1).Value = '06:52:34'
$ws3.Cells.Item(13,2).Value = '07:00'
$ws3.Cells.Item(13,3).Value = '215B_LP-P MOR-1 Y 57'
$ws3.Cells.Item(13,4).Value = 8
$ws3.Cells.Item(13,5).Value = 'L6173'

$ws3.Cells.Item(14,1).Value = '06:52:34'
$ws3.Cells.Item(14,2).Value = '07:35'
$ws3.Cells.Item(14,3).Value = '215A_LA PLATA'
$ws3.Cells.Item(14,4).Value = 43
$ws3.Cells.Item(14,5).Value = 'L6173'

$ws3.Cells.Item(15,1).Value = '06:52:34'
$ws3.Cells.Item(15,2).Value = '08:06'
$ws3.Cells.Item(15,3).Value = '215C_LA PLATA'
$ws3.Cells.Item(15,4).Value = 74
$ws3.Cells.Item(15,5).Value = 'L6203'

$ws3.Cells.Item(16,1).Value = '06:52:34'
$ws3.Cells.Item(16,2).Value = '08:33'
$ws3.Cells.Item(16,3).Value = '215A_LA PLATA'
$ws3.Cells.Item(16,4).Value = 101
$ws3.Cells.Item(16,5).Value = 'L6173'

Write-Output "Edit complete"